# Weekly price-log update: a new daily record was added for
# "Feria Lagunitas de Puerto Montt - Coliflor", inserted as row 274,
# pushing every following record down by one row (274-337 -> 275-338).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 274 (shifts old rows 274..337 down to 275..338,
# and grows the used range to A1:R338).
$ws.Rows.Item(274).EntireRow.Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A274").Value = 4
$ws.Range("B274").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C274").Value = "Los Lagos"
$ws.Range("D274").Value = 44722
$ws.Range("E274").Value = 10
$ws.Range("F274").Value = 100112008
$ws.Range("G274").Value = "Coliflor"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Primera"
$ws.Range("J274").Value = 1000
$ws.Range("K274").Value = 1600
$ws.Range("L274").Value = 1700
$ws.Range("M274").Value = 1650
$ws.Range("N274").Value = "$/unidad"
$ws.Range("O274").Value = "Región Metropolitana"
$ws.Range("P274").Value = 1650
$ws.Range("Q274").Value = 1
$ws.Range("R274").Value = "Hortaliza"
